$wb = $excel.ActiveWorkbook

# --- Rename "BCS" -> "BCS-BCS" and add the new "BCS-DoSfCS" sheet right after it ---
$bcs = $wb.Worksheets.Item("BCS")
$bcs.Name = "BCS-BCS"

$dosfcs = $wb.Worksheets.Add($null, $bcs)
$dosfcs.Name = "BCS-DoSfCS"
$dosfcs.Tab.Color = 8210719

# Header row: "Years" label + the 2021-2050 year series
$dosfcs.Range("A1").Value = "Years"
$col = 2
foreach ($y in 2021..2050) {
    $dosfcs.Cells.Item(1, $col).Value = $y
    $col++
}

# Duration row: 12-year 45Q credit duration, repeated via formula reference to $B$2
$dosfcs.Range("A2").Value = "Duration"
$dosfcs.Range("B2").Value = 12
$dosfcs.Range("C2").Formula = "=`$B`$2"
$dosfcs.Range("D2").Formula = "=`$B`$2"
$dosfcs.Range("E2:AE2").Formula = "=`$B`$2"

$dosfcs.Range("C32").Select() | Out-Null

# --- About sheet: document the 45Q credit duration used above ---
$about = $wb.Worksheets.Item("About")
$about.Range("A13").Value = "45Q Duration"
$about.Range("B13").Value = "12 years"

# --- BCS-BCS: cost-per-unit new elec output now keys off About!$B$11 instead of the ---
# --- Electricity Calculations dispatch-cost figure, and loses its currency number format ---
$cols = @("D","E","F","G","H","I","J","K","L","M")
foreach ($c in $cols) {
    $cell = $bcs.Range("$c`2")
    $cell.Formula = "=About!`$B`$11*About!`$A`$9"
    $cell.Style = "Normal"
}
$bcs.Rows.Item(2).AutoFit()
$bcs.Rows.Item(3).AutoFit()

$about.Activate() | Out-Null
$about.Range("A14").Select() | Out-Null
